# "Added Week 15 simulations" — update Rushing and Receiving sheets with the
# new week's data: a few existing players' cumulative totals change, a new
# D.Swift row is added to both sheets, a new C.Reynolds row is added to
# Receiving, and new trailing rows (J.Jefferson on Rushing, S.Zylstra on
# Receiving) are appended.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Rushing sheet
# ---------------------------------------------------------------------
$rushing = $wb.Worksheets.Item("Rushing")

# Existing rows 2-9: update name/values in place (no row shifting needed —
# final values below already account for the new week's data).
$rushing.Range("A2").Value = 0
$rushing.Range("B2").Value = "J.Goff"
$rushing.Range("C2").Value = 7
$rushing.Range("D2").Value = 3
$rushing.Range("E2").Value = 3
$rushing.Range("F2").Value = 2

$rushing.Range("A3").Value = 1
$rushing.Range("B3").Value = "D.Swift"
$rushing.Range("C3").Value = 74
$rushing.Range("D3").Value = 45
$rushing.Range("E3").Value = 21
$rushing.Range("F3").Value = 13

$rushing.Range("A4").Value = 2
$rushing.Range("B4").Value = "J.Williams"
$rushing.Range("C4").Value = 54
$rushing.Range("D4").Value = 27
$rushing.Range("E4").Value = 13
$rushing.Range("F4").Value = 11

$rushing.Range("A5").Value = 3
$rushing.Range("B5").Value = "J.Jefferson"
$rushing.Range("C5").Value = 1
$rushing.Range("D5").Value = 0
$rushing.Range("E5").Value = 1
$rushing.Range("F5").Value = 1

$rushing.Range("A6").Value = 4
$rushing.Range("B6").Value = "G.Igwebuike"
$rushing.Range("C6").Value = 3
$rushing.Range("D6").Value = 1
$rushing.Range("E6").Value = 0
$rushing.Range("F6").Value = 0

$rushing.Range("A7").Value = 5
$rushing.Range("B7").Value = "J.Cabinda"
$rushing.Range("C7").Value = 0
$rushing.Range("D7").Value = 1
$rushing.Range("E7").Value = 2
$rushing.Range("F7").Value = 0

$rushing.Range("A8").Value = 6
$rushing.Range("B8").Value = "K.Raymond"
$rushing.Range("C8").Value = 0
$rushing.Range("D8").Value = 2
$rushing.Range("E8").Value = 0
$rushing.Range("F8").Value = 0

$rushing.Range("A9").Value = 7
$rushing.Range("B9").Value = "J.Jefferson"
$rushing.Range("C9").Value = 1
$rushing.Range("D9").Value = 0
$rushing.Range("E9").Value = 1
$rushing.Range("F9").Value = 1

# New trailing row 10 — copy the number-format/border/bold style of column A
# (style index 1 in the original file) from A9 onto A10 before setting values.
$rushing.Range("A9").Copy()
$rushing.Range("A10").PasteSpecial(-4122)
$rushing.Range("A10").Value = 8
$rushing.Range("B10").Value = "J.Jefferson"
$rushing.Range("C10").Value = 1
$rushing.Range("D10").Value = 0
$rushing.Range("E10").Value = 1
$rushing.Range("F10").Value = 1

# ---------------------------------------------------------------------
# Receiving sheet
# ---------------------------------------------------------------------
$receiving = $wb.Worksheets.Item("Receiving")

# New lead row for D.Swift's week — overwrite old row 2 (J.Williams) values,
# then rewrite every following row with the shifted-down content.
$receiving.Range("A2").Value = 0
$receiving.Range("B2").Value = "D.Swift"
$receiving.Range("C2").Value = 67
$receiving.Range("D2").Value = 53
$receiving.Range("E2").Value = 3
$receiving.Range("F2").Value = 3
$receiving.Range("G2").Value = 7
$receiving.Range("H2").Value = 4

$receiving.Range("A3").Value = 1
$receiving.Range("B3").Value = "J.Williams"
$receiving.Range("C3").Value = 23
$receiving.Range("D3").Value = 21
$receiving.Range("E3").Value = 0
$receiving.Range("F3").Value = 0
$receiving.Range("G3").Value = 0
$receiving.Range("H3").Value = 0

$receiving.Range("A4").Value = 2
$receiving.Range("B4").Value = "G.Igwebuike"
$receiving.Range("C4").Value = 7
$receiving.Range("D4").Value = 5
$receiving.Range("E4").Value = 0
$receiving.Range("F4").Value = 0
$receiving.Range("G4").Value = 0
$receiving.Range("H4").Value = 0

$receiving.Range("A5").Value = 3
$receiving.Range("B5").Value = "J.Cabinda"
$receiving.Range("C5").Value = 2
$receiving.Range("D5").Value = 0
$receiving.Range("E5").Value = 0
$receiving.Range("F5").Value = 0
$receiving.Range("G5").Value = 0
$receiving.Range("H5").Value = 0

# New row for C.Reynolds
$receiving.Range("A6").Value = 4
$receiving.Range("B6").Value = "C.Reynolds"
$receiving.Range("C6").Value = 2
$receiving.Range("D6").Value = 2
$receiving.Range("E6").Value = 0
$receiving.Range("F6").Value = 0
$receiving.Range("G6").Value = 0
$receiving.Range("H6").Value = 0

$receiving.Range("A7").Value = 5
$receiving.Range("B7").Value = "T.Williams"
$receiving.Range("C7").Value = 3
$receiving.Range("D7").Value = 2
$receiving.Range("E7").Value = 0
$receiving.Range("F7").Value = 0
$receiving.Range("G7").Value = 1
$receiving.Range("H7").Value = 1

$receiving.Range("A8").Value = 6
$receiving.Range("B8").Value = "A.St. Brown"
$receiving.Range("C8").Value = 38
$receiving.Range("D8").Value = 31
$receiving.Range("E8").Value = 5
$receiving.Range("F8").Value = 2
$receiving.Range("G8").Value = 4
$receiving.Range("H8").Value = 2

$receiving.Range("A9").Value = 7
$receiving.Range("B9").Value = "K.Raymond"
$receiving.Range("C9").Value = 41
$receiving.Range("D9").Value = 31
$receiving.Range("E9").Value = 16
$receiving.Range("F9").Value = 7
$receiving.Range("G9").Value = 4
$receiving.Range("H9").Value = 3

$receiving.Range("A10").Value = 8
$receiving.Range("B10").Value = "Q.Cephus"
$receiving.Range("C10").Value = 13
$receiving.Range("D10").Value = 10
$receiving.Range("E10").Value = 10
$receiving.Range("F10").Value = 5
$receiving.Range("G10").Value = 3
$receiving.Range("H10").Value = 3

$receiving.Range("A11").Value = 9
$receiving.Range("B11").Value = "K.Hodge"
$receiving.Range("C11").Value = 11
$receiving.Range("D11").Value = 7
$receiving.Range("E11").Value = 5
$receiving.Range("F11").Value = 0
$receiving.Range("G11").Value = 2
$receiving.Range("H11").Value = 0

$receiving.Range("A12").Value = 10
$receiving.Range("B12").Value = "T.Benson"
$receiving.Range("C12").Value = 13
$receiving.Range("D12").Value = 7
$receiving.Range("E12").Value = 5
$receiving.Range("F12").Value = 1
$receiving.Range("G12").Value = 1
$receiving.Range("H12").Value = 0

$receiving.Range("A13").Value = 11
$receiving.Range("B13").Value = "T.Kennedy"
$receiving.Range("C13").Value = 2
$receiving.Range("D13").Value = 2
$receiving.Range("E13").Value = 1
$receiving.Range("F13").Value = 0
$receiving.Range("G13").Value = 0
$receiving.Range("H13").Value = 0

$receiving.Range("A14").Value = 12
$receiving.Range("B14").Value = "J.Jefferson"
$receiving.Range("C14").Value = 4
$receiving.Range("D14").Value = 4
$receiving.Range("E14").Value = 0
$receiving.Range("F14").Value = 0
$receiving.Range("G14").Value = 0
$receiving.Range("H14").Value = 0

$receiving.Range("A15").Value = 13
$receiving.Range("B15").Value = "B.Wright"
$receiving.Range("C15").Value = 10
$receiving.Range("D15").Value = 6
$receiving.Range("E15").Value = 1
$receiving.Range("F15").Value = 1
$receiving.Range("G15").Value = 2
$receiving.Range("H15").Value = 0

$receiving.Range("A16").Value = 14
$receiving.Range("B16").Value = "J.Reynolds"
$receiving.Range("C16").Value = 10
$receiving.Range("D16").Value = 7
$receiving.Range("E16").Value = 8
$receiving.Range("F16").Value = 4
$receiving.Range("G16").Value = 1
$receiving.Range("H16").Value = 0

$receiving.Range("A17").Value = 15
$receiving.Range("B17").Value = "T.Hockenson"
$receiving.Range("C17").Value = 71
$receiving.Range("D17").Value = 53
$receiving.Range("E17").Value = 13
$receiving.Range("F17").Value = 8
$receiving.Range("G17").Value = 10
$receiving.Range("H17").Value = 7

$receiving.Range("A18").Value = 16
$receiving.Range("B18").Value = "D.Fells"
$receiving.Range("C18").Value = 4
$receiving.Range("D18").Value = 3
$receiving.Range("E18").Value = 1
$receiving.Range("F18").Value = 1
$receiving.Range("G18").Value = 0
$receiving.Range("H18").Value = 0

# New trailing row 19 — copy column-A style from row 17 onto the new row.
$receiving.Range("A17").Copy()
$receiving.Range("A19").PasteSpecial(-4122)
$receiving.Range("A19").Value = 17
$receiving.Range("B19").Value = "S.Zylstra"
$receiving.Range("C19").Value = 5
$receiving.Range("D19").Value = 2
$receiving.Range("E19").Value = 1
$receiving.Range("F19").Value = 0
$receiving.Range("G19").Value = 2
$receiving.Range("H19").Value = 0

# ---------------------------------------------------------------------
# Selection / active-sheet state (matches the saved view in the workbook)
# ---------------------------------------------------------------------
$rushing.Range("F12").Select()
$receiving.Range("I2").Select()
$receiving.Activate()
